# Fix sample numbers 001-006 causing false positive data columns in PCA plots.
# Updates the "UnitMass" column (C) values on the active worksheet for both
# the "+ loading" table (rows 2-21) and the "- loading" table (rows 23-42).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "C2" = 23
    "C3" = 43
    "C4" = 58
    "C5" = 85
    "C6" = 70
    "C7" = 111
    "C8" = 31
    "C9" = 98
    "C10" = 55
    "C11" = 19
    "C12" = 57
    "C13" = 71
    "C14" = 99
    "C15" = 53
    "C16" = 60
    "C17" = 110
    "C18" = 66
    "C19" = 86
    "C20" = 138
    "C21" = 72
    "C23" = 39
    "C24" = 46
    "C25" = 27
    "C26" = 32
    "C27" = 91
    "C28" = 125
    "C29" = 83
    "C30" = 40
    "C31" = 28
    "C32" = 88
    "C33" = 175
    "C34" = 56
    "C35" = 231
    "C36" = 7
    "C37" = 24
    "C38" = 106
    "C39" = 115
    "C40" = 159
    "C41" = 90
    "C42" = 1
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
